$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$used.Replace("D64", "D69")
$used.Replace("D80", "D86")
$used.Replace("D51", "D55")
$used.Replace("S30", "S31")
